$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.753.27"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.288.48"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.65"
$ws.Range("E5").Value = "  +14.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.32"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.69"
$ws.Range("E10").Value = "  +5.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.06"
$ws.Range("E12").Value = "  +14.88%  "
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.75"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").Value = "2.629.58"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "2.283.61"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "43.655.44"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.73"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.16"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -5.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.91"
$ws.Range("E23").Value = "  +8.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.32"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("E25").Value = "  +5.74%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +4.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.82"
$ws.Range("E28").Value = "  +4.91%  "
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.87"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.53"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0926"
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("E34").Value = "  +5.30%  "
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.68"
$ws.Range("E36").Value = "  +7.22%  "
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("E39").Value = "  +14.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.244"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.68"
$ws.Range("E41").Value = "  +11.20%  "
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.39"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.74"
$ws.Range("E43").Value = "  +11.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.31"
$ws.Range("E44").Value = "  +22.72%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.46"
$ws.Range("E48").Value = "  +5.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0993"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("E51").Value = "  +5.13%  "
